$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph (the footer/nav line that, together
# with the blank paragraph before it and the copyright paragraph after it,
# needs to be removed from the end of the "Requisitos" section).
$findRange = $d.Content
$findRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$verParaIndex = $findRange.Paragraphs(1).Index

# Remove the blank paragraph right before it, the "Ver no Jupiter..." paragraph
# itself, and the "(c) 2020 ..." copyright paragraph right after it, leaving
# the trailing blank paragraph and page-break paragraph untouched.
$startPara = $d.Paragraphs($verParaIndex - 1)
$endPara = $d.Paragraphs($verParaIndex + 1)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
